# Updates the cryptos list (Price / Volume(1h) columns, and the
# Mantle / BabyDogeCoin row swap) to match the latest scrape.
# Every written cell is forced to Text (leading apostrophe) and then
# restyled back to "Normal" so no stray number-format/quote-prefix
# style sticks around -- the source file keeps these as plain inline
# strings with the default (unstyled) cell format.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.135.44"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +0.17%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.838.37"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +0.22%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.45%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'243.16"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -0.79%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'0.6175"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -2.71%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'1.004"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +0.43%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.07478"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -1.12%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.2934"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -0.58%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'23.14"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -0.49%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.07704"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -0.24%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'1.831.92"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -0.15%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'5.004"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -0.16%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.6732"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -0.17%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'82.83"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -0.56%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.000009150"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -4.81%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'5.922"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -3.05%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'29.096.54"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -0.01%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'2.078.14"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -0.20%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'233.51"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +2.53%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'12.68"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +0.40%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'1.004"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +0.52%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'7.194"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +0.10%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'  +0.49%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'159.36"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -0.61%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.1406"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -1.93%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'8.500"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -0.68%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'17.87"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -0.57%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'1.502"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -0.07%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'4.161"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +0.03%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'4.113"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +0.76%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'0.05514"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +0.75%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'1.210"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +0.63%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'1.836"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -1.63%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'0.7383"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -1.32%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'1.143"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +0.03%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'2.667"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +0.27%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'2.774"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +0.50%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.01783"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -0.32%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'1.213.46"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -2.87%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'6.456"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -2.89%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.8924"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -1.63%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'1.003"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +0.29%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'102.04"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +0.25%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'1.980.21"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -0.20%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'65.54"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +0.26%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("B47").Value = "'BabyDogeCoin"
$ws.Range("B47").Style = "Normal"
$ws.Range("C47").Value = "'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("C47").Style = "Normal"
$ws.Range("D47").Value = "'0.00000000121"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -0.60%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("B48").Value = "'Mantle"
$ws.Range("B48").Style = "Normal"
$ws.Range("C48").Value = "'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("C48").Style = "Normal"
$ws.Range("D48").Value = "'0.5099"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -0.15%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.4079"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -0.12%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'9.167"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +1.59%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.05827"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +0.57%  "
$ws.Range("E51").Style = "Normal"
